$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AZ (column 52), pushing "nom" (was AZ) to BA
# and "url_produit" (was BA) to BB, matching the site's new price-history
# snapshot column added ahead of the name/url columns.
$ws.Columns.Item(52).Insert()

# New snapshot header timestamp
$ws.Range("AZ1").Value = "2026-01-29 23:16:25"

# Populate the new snapshot column (AZ) with each product's latest known
# price carried forward from the previous snapshot (AY), for rows that
# still had a numeric price there. Rows whose AY was already blank
# (discontinued/out of stock) stay blank in AZ too - the freshly inserted
# column cells are already empty by default, so nothing to do there.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ayCell = $ws.Cells.Item($r, 51)
    $ayVal = $ayCell.Value()
    if ($ayVal -ne $null -and $ayVal -ne "") {
        $ws.Cells.Item($r, 52).Value = $ayVal
    }
}
